# Automatic update of files.
# Bump the "Förändrad" (Changed) date in column C for every data row
# (rows 2-52) from 45188 (2023-09-19) to 45189 (2023-09-20), keeping the
# existing date formatting intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 52; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
